# Applies the "Add files via upload" commit to the 16-9-1 metadata workbook.
#
# Summary of the change (from the OOXML diff):
#   * The author refreshed the contact / organisation details for the
#     16.9.1 indicator sheet:
#       - B4  (Индикатор)   -> drop the trailing period after "16.9.1"
#       - B6  (Организация) -> "Отдел" -> "Управление" статистики домашних хозяйств
#       - B7  (Контактное лицо) -> unchanged text, kept as-is
#       - B8  (Электронная почта) -> y.kalymbetova@ -> yryskan.kalymbetova@
#       - B9  (Телефон) -> unchanged text, kept as-is
#       - B10 (Сайт) -> www.stat.kg -> www.stat.gov.kg
#   * The workbook window was left maximised (0,0 / 28800x11835) and the
#     surviving selection anchor is B6.
#
# NOTE: the underlying B7/B9/... values/text did not actually change in the
# diff (only the shared-string slot they point at did, because four stale
# strings were dropped from the table) so we leave those cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 2. Maximise / reposition the workbook window (best effort; cosmetic) ---
$excel.WindowState = -4137   # xlMaximized
$excel.Left   = 0
$excel.Top    = 0
$excel.Width  = 28800
$excel.Height = 11835

$win = $excel.ActiveWindow
$win.WindowState = -4137   # xlMaximized
$win.Left   = 0
$win.Top    = 0
$win.Width  = 28800
$win.Height = 11835

# --- 1. Refresh the indicator / organisation / contact text -----------------
$ws.Range("B4").Value  = "16.9.1 Доля детей в возрасте до пяти лет, рождение которых было зарегистрировано в гражданских органах, в разбивке по возрасту"
$ws.Range("B6").Value  = "Национальный статистический комитет Кыргызской Республики (Управление статистики домашних хозяйств)"
$ws.Range("B8").Value  = "yryskan.kalymbetova@gmail.com "
$ws.Range("B10").Value = "www.stat.gov.kg"

# --- 3. Move the visible selection to B6 ------------------------------------
$ws.Range("B6").Select()
